$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (series C)
$ws.Range("J2").Value = 241.1142161815241
$ws.Range("K2").Value = 283.1250585464616
$ws.Range("L2").Value = 205.3370533265266
$ws.Range("M2").Value = 4.518851520511773
$ws.Range("N2").Value = 8.149026565380112
$ws.Range("O2").Value = 2.505823106674223
$ws.Range("P2").Value = 124.958401618425
$ws.Range("T2").Value = 11.35491950195709
$ws.Range("U2").Value = 90.96286852303754
$ws.Range("Z2").Value = "241`n(205, 283)"
$ws.Range("AA2").Value = "4.5`n(2.5, 8.1)"

# Row 3 (series D)
$ws.Range("J3").Value = 134.1449701011463
$ws.Range("K3").Value = 144.5290465001201
$ws.Range("L3").Value = 124.5069654799285
$ws.Range("M3").Value = 8.768954184498069
$ws.Range("N3").Value = 14.44228397562113
$ws.Range("O3").Value = 5.324265719994549
$ws.Range("P3").Value = 95.60297344464202
$ws.Range("T3").Value = 24.51663412987911
$ws.Range("U3").Value = 133.2753708746945
$ws.Range("Z3").Value = "134`n(125, 145)"
$ws.Range("AA3").Value = "8.8`n(5.3, 14.4)"

# Row 4 (series B)
$ws.Range("J4").Value = 146.8346489901402
$ws.Range("K4").Value = 170.1596387263386
$ws.Range("L4").Value = 126.7069811938922
$ws.Range("M4").Value = 4.42582981773189
$ws.Range("N4").Value = 7.335649771765482
$ws.Range("O4").Value = 2.670243289274487
$ws.Range("P4").Value = 75.05359092245527
$ws.Range("T4").Value = 110.0847982349058
$ws.Range("U4").Value = 269.5831474759362
$ws.Range("Z4").Value = "147`n(127, 170)"
$ws.Range("AA4").Value = "4.4`n(2.7, 7.3)"

# Row 5 (series E)
$ws.Range("J5").Value = 80.66718380179394
$ws.Range("K5").Value = 88.31352716578164
$ws.Range("L5").Value = 73.68287454193896
$ws.Range("M5").Value = 11.42761107958443
$ws.Range("N5").Value = 24.74932316464394
$ws.Range("O5").Value = 5.276519851371034
$ws.Range("P5").Value = 62.20391396429272
$ws.Range("T5").Value = 181.5729904224387
$ws.Range("U5").Value = 302.3052074333061
$ws.Range("Z5").Value = "81`n(74, 88)"
$ws.Range("AA5").Value = "11.4`n(5.3, 24.7)"

# Row 8 (series F)
$ws.Range("J8").Value = 64.03564641645039
$ws.Range("K8").Value = 74.43443826719667
$ws.Range("L8").Value = 55.08960781369634
$ws.Range("M8").Value = 6.826103861777176
$ws.Range("N8").Value = 15.52715750065982
$ws.Range("O8").Value = 3.000915906841881
$ws.Range("P8").Value = 41.44291078397361
$ws.Range("T8").Value = 899.9401649186126
$ws.Range("U8").Value = 655.765974403018
$ws.Range("Z8").Value = "64`n(55, 74)"
$ws.Range("AA8").Value = "6.8`n(3.0, 15.5)"
